# Guinea GDP per Capita workbook update: refresh the GDP-per-capita series
# with revised figures and extend the "Data" sheet through 2016.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- 1. Add new rows 61-68 (years 2009-2016) with their fixed columns ---
$newRows = @(
    @{Row=61; Year=2009},
    @{Row=62; Year=2010},
    @{Row=63; Year=2011},
    @{Row=64; Year=2012},
    @{Row=65; Year=2013},
    @{Row=66; Year=2014},
    @{Row=67; Year=2015},
    @{Row=68; Year=2016}
)
foreach ($nr in $newRows) {
    $row = $nr.Row
    $ws.Range("A" + $row).Value = 324
    $ws.Range("B" + $row).Value = "Guinea"
    $ws.Range("C" + $row).Value = "GDP per Capita"
    $ws.Range("D" + $row).Value = $nr.Year
}

# --- 2. Refresh the "Data" value (column E) for every row 2-68. ---
# The source values are numeric-looking strings that must stay stored as text
# (shared strings), exactly like the rest of the column. Assigning them straight
# to .Value would let Excel auto-coerce them to numbers, so instead each value is
# written as a literal-string formula in a scratch column, then copied and pasted
# back into column E with Paste Special > Values - this keeps the cells text-typed
# without touching number formats/styles.
$eValues = @(
    @{Row=2; Val="520"},
    @{Row=3; Val="542"},
    @{Row=4; Val="550"},
    @{Row=5; Val="563"},
    @{Row=6; Val="587"},
    @{Row=7; Val="598"},
    @{Row=8; Val="610"},
    @{Row=9; Val="622"},
    @{Row=10; Val="630"},
    @{Row=11; Val="655"},
    @{Row=12; Val="671"},
    @{Row=13; Val="703"},
    @{Row=14; Val="741"},
    @{Row=15; Val="689"},
    @{Row=16; Val="720"},
    @{Row=17; Val="756"},
    @{Row=18; Val="759"},
    @{Row=19; Val="767"},
    @{Row=20; Val="776"},
    @{Row=21; Val="784"},
    @{Row=22; Val="795"},
    @{Row=23; Val="822"},
    @{Row=24; Val="826"},
    @{Row=25; Val="821"},
    @{Row=26; Val="848"},
    @{Row=27; Val="854"},
    @{Row=28; Val="905"},
    @{Row=29; Val="858"},
    @{Row=30; Val="897"},
    @{Row=31; Val="869"},
    @{Row=32; Val="872"},
    @{Row=33; Val="859"},
    @{Row=34; Val="858"},
    @{Row=35; Val="850"},
    @{Row=36; Val="824"},
    @{Row=37; Val="821"},
    @{Row=38; Val="819"},
    @{Row=39; Val="826"},
    @{Row=40; Val="856"},
    @{Row=41; Val="869"},
    @{Row=42; Val="861"},
    @{Row=43; Val="827.006899123684"},
    @{Row=44; Val="833.226762332466"},
    @{Row=45; Val="862.731285361005"},
    @{Row=46; Val="882.494426799569"},
    @{Row=47; Val="906.264389671341"},
    @{Row=48; Val="927.12947969862"},
    @{Row=49; Val="976.496625957152"},
    @{Row=50; Val="1010.79190263154"},
    @{Row=51; Val="1031.38254533899"},
    @{Row=52; Val="1046.56706272681"},
    @{Row=53; Val="1090.65199691792"},
    @{Row=54; Val="1148.54467847906"},
    @{Row=55; Val="1146.3352913077"},
    @{Row=56; Val="1159.82503467527"},
    @{Row=57; Val="1186.60479343225"},
    @{Row=58; Val="1188.48457591907"},
    @{Row=59; Val="1253.40771871762"},
    @{Row=60; Val="1290.4037004352"},
    @{Row=61; Val="1255.24185529761"},
    @{Row=62; Val="1291.69622840926"},
    @{Row=63; Val="1347"},
    @{Row=64; Val="1397"},
    @{Row=65; Val="1421"},
    @{Row=66; Val="1399"},
    @{Row=67; Val="1364"},
    @{Row=68; Val="1417"}
)
foreach ($e in $eValues) {
    $ws.Range("ZZ" + $e.Row).Formula = '="' + $e.Val + '"'
}

$ws.Range("ZZ2:ZZ68").Copy() | Out-Null
$ws.Range("E2:E68").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range("ZZ2:ZZ68").Clear()

Write-Host "Guinea GDP per Capita data refreshed through 2016."